$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the slightly refined timestamp value on the existing last row (63)
$ws.Cells.Item(63, 1).Value = 44376.76795573263

# Append the newly retrieved data row (64)
$ws.Cells.Item(64, 1).Value = 44377.76845540029
$ws.Cells.Item(64, 2).Value = 78440
$ws.Cells.Item(64, 3).Value = 66002
$ws.Cells.Item(64, 4).Value = 3651
$ws.Cells.Item(64, 5).Value = 2140
$ws.Cells.Item(64, 6).Value = 1525
$ws.Cells.Item(64, 7).Value = 20889
$ws.Cells.Item(64, 8).Value = 1628
$ws.Cells.Item(64, 9).Value = 878
$ws.Cells.Item(64, 10).Value = 207

# Ensure the new date cell uses the same number format as the rest of column A
$ws.Cells.Item(64, 1).NumberFormat = "yyyy-mm-dd HH:mm:ss UTC"
